$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 3100
$ws.Range("I2").Value = 5250
$ws.Range("J2").Value = 950
$ws.Range("K2").Value = 5250
$ws.Range("L2").Value = 950
$ws.Range("M2").Value = -5137
$ws.Range("N2").Value = -1176
$ws.Range("H15").Value = 2124.5454
$ws.Range("I15").Value = 2124.5454
$ws.Range("K15").Value = 6373.6362
$ws.Range("M15").Value = -6204.6362
$ws.Range("H64").Value = 4490
$ws.Range("J64").Value = 3990
$ws.Range("L64").Value = 3990
$ws.Range("N64").Value = -4486
$ws.Range("H67").Value = 4490
$ws.Range("J67").Value = 3990
$ws.Range("L67").Value = 3990
$ws.Range("N67").Value = -5706
$ws.Range("H76").Value = 2746.4546
$ws.Range("I76").Value = 2690.85
$ws.Range("K76").Value = 2690.85
$ws.Range("M76").Value = -2375.85
$ws.Range("H79").Value = 2746.4546
$ws.Range("I79").Value = 2690.85
$ws.Range("K79").Value = 2690.85
$ws.Range("M79").Value = -1598.85

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5001237.5
$ws.Range("I2").Value = 8929473
$ws.Range("J2").Value = 1664.5
$ws.Range("K2").Value = 8929473
$ws.Range("L2").Value = 1664.5
$ws.Range("M2").Value = -8929360
$ws.Range("N2").Value = -1890.5
$ws.Range("H32").Value = 3632.9
$ws.Range("I32").Value = 3632.9
$ws.Range("K32").Value = 3632.9
$ws.Range("M32").Value = -3345.9
$ws.Range("H63").Value = 3099
$ws.Range("I63").Value = 1748.3
$ws.Range("J63").Value = 5028.5713
$ws.Range("K63").Value = 1748.3
$ws.Range("L63").Value = 5028.5713
$ws.Range("M63").Value = -1062.3
$ws.Range("N63").Value = -6400.5713
$ws.Range("H66").Value = 3099
$ws.Range("I66").Value = 1748.3
$ws.Range("J66").Value = 5028.5713
$ws.Range("K66").Value = 8741.5
$ws.Range("L66").Value = 25142.8565
$ws.Range("M66").Value = -5309.5
$ws.Range("N66").Value = -32006.8565
$ws.Range("H116").Value = 5001237.5
$ws.Range("I116").Value = 8929473
$ws.Range("J116").Value = 1664.5
$ws.Range("K116").Value = 8929473
$ws.Range("L116").Value = 1664.5
$ws.Range("M116").Value = -8927179
$ws.Range("N116").Value = -6252.5

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5001237.5
$ws.Range("I3").Value = 8929473
$ws.Range("J3").Value = 1664.5
$ws.Range("K3").Value = 8929473
$ws.Range("L3").Value = 1664.5
$ws.Range("M3").Value = -8929359
$ws.Range("N3").Value = -1892.5
$ws.Range("H105").Value = 1607.8889
$ws.Range("I105").Value = 1465.7142
$ws.Range("J105").Value = 2105.5
$ws.Range("K105").Value = 1465.7142
$ws.Range("L105").Value = 2105.5
$ws.Range("M105").Value = 281.2858000000001
$ws.Range("N105").Value = -5599.5
$ws.Range("H134").Value = 3137
$ws.Range("I134").Value = 2979.5557
$ws.Range("J134").Value = 3924.2222
$ws.Range("K134").Value = 8938.667099999999
$ws.Range("L134").Value = 11772.6666
$ws.Range("M134").Value = -6403.667099999999
$ws.Range("N134").Value = -16842.6666

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3510.6843
$ws.Range("I62").Value = 2480.8333
$ws.Range("J62").Value = 3986
$ws.Range("K62").Value = 2480.8333
$ws.Range("L62").Value = 3986
$ws.Range("M62").Value = -1856.8333
$ws.Range("N62").Value = -5234
$ws.Range("H65").Value = 3510.6843
$ws.Range("I65").Value = 2480.8333
$ws.Range("J65").Value = 3986
$ws.Range("K65").Value = 12404.1665
$ws.Range("L65").Value = 19930
$ws.Range("M65").Value = -9284.166499999999
$ws.Range("N65").Value = -26170
$ws.Range("H134").Value = 13891485
$ws.Range("I134").Value = 19232752
$ws.Range("J134").Value = 4191.3
$ws.Range("K134").Value = 57698256
$ws.Range("L134").Value = 12573.9
$ws.Range("M134").Value = -57695721
$ws.Range("N134").Value = -17643.9

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1779.9
$ws.Range("I97").Value = 433.33334
$ws.Range("J97").Value = 2357
$ws.Range("K97").Value = 1300.00002
$ws.Range("L97").Value = 7071
$ws.Range("M97").Value = -804.0000199999999
$ws.Range("N97").Value = -8063
$ws.Range("H131").Value = 821.24
$ws.Range("I131").Value = 320.92307
$ws.Range("J131").Value = 896
$ws.Range("K131").Value = 962.7692099999999
$ws.Range("L131").Value = 2688
$ws.Range("M131").Value = 4077.23079
$ws.Range("N131").Value = -12768

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3567.36
$ws.Range("I80").Value = 3106.8572
$ws.Range("J80").Value = 4153.4546
$ws.Range("K80").Value = 3106.8572
$ws.Range("L80").Value = 4153.4546
$ws.Range("M80").Value = -2108.8572
$ws.Range("N80").Value = -6149.4546
$ws.Range("H83").Value = 3567.36
$ws.Range("I83").Value = 3106.8572
$ws.Range("J83").Value = 4153.4546
$ws.Range("K83").Value = 15534.286
$ws.Range("L83").Value = 20767.273
$ws.Range("M83").Value = -10542.286
$ws.Range("N83").Value = -30751.273
$ws.Range("H132").Value = 3669.1177
$ws.Range("I132").Value = 3640.2122
$ws.Range("K132").Value = 10920.6366
$ws.Range("M132").Value = -8390.6366
$ws.Range("H136").Value = 8620.233
$ws.Range("J136").Value = 8620.233
$ws.Range("L136").Value = 25860.699
$ws.Range("N136").Value = -30960.699

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2761.0588
$ws.Range("I82").Value = 2161.6667
$ws.Range("J82").Value = 4199.6
$ws.Range("K82").Value = 2161.6667
$ws.Range("L82").Value = 4199.6
$ws.Range("M82").Value = -1800.6667
$ws.Range("N82").Value = -4921.6
$ws.Range("H85").Value = 2761.0588
$ws.Range("I85").Value = 2161.6667
$ws.Range("J85").Value = 4199.6
$ws.Range("K85").Value = 2161.6667
$ws.Range("L85").Value = 4199.6
$ws.Range("M85").Value = -913.6667000000002
$ws.Range("N85").Value = -6695.6
$ws.Range("H132").Value = 2713.2
$ws.Range("I132").Value = 1840.8096
$ws.Range("J132").Value = 4021.7856
$ws.Range("K132").Value = 5522.4288
$ws.Range("L132").Value = 12065.3568
$ws.Range("M132").Value = -2992.4288
$ws.Range("N132").Value = -17125.3568

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 38967.5
$ws.Range("J137").Value = 38967.5
$ws.Range("L137").Value = 38967.5
$ws.Range("N137").Value = -49167.5
